$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# Update the "Jeff, Keith, Phil, Andy, Omar, Brad" roster text -> swap Omar for Andrei
$ws.Range("A10").Value = "Jeff, Keith, Phil, Andy, Brad, Andrei"

# Remove the obsolete "House Hunting  Backlog" row entirely (row 14), shifting
# everything below it up by one row
$ws.Rows(14).Delete()

# Update the active selection to match the author's final cursor position
$ws.Range("A11").Select()
